# Sync attendance_reports, modules_schedules, and assets from main repo - 2026-01-26 01:22:45
# Update the "Recorded By" column (G) values so that "System"/"system" is
# listed last among the recorder names instead of first.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @{
    "G2"   = "system, backup@backdoor.com, System"
    "G3"   = "dnasr281@gmail.com, System"
    "G5"   = "backup@backdoor.com, System"
    "G6"   = "dnasr281@gmail.com, System"
    "G7"   = "admin@admin.com, System"
    "G8"   = "backup@backdoor.com, System"
    "G28"  = "system, backup@backdoor.com, System"
    "G29"  = "dnasr281@gmail.com, System"
    "G31"  = "backup@backdoor.com, System"
    "G32"  = "dnasr281@gmail.com, System"
    "G33"  = "admin@admin.com, System"
    "G34"  = "backup@backdoor.com, System"
    "G54"  = "system, backup@backdoor.com, System"
    "G55"  = "dnasr281@gmail.com, System"
    "G57"  = "backup@backdoor.com, System"
    "G58"  = "dnasr281@gmail.com, System"
    "G59"  = "admin@admin.com, System"
    "G60"  = "backup@backdoor.com, System"
    "G80"  = "backup@backdoor.com, System"
    "G81"  = "backup@backdoor.com, System"
    "G82"  = "backup@backdoor.com, System"
    "G87"  = "admin@admin.com, dnasr281@gmail.com"
    "G106" = "backup@backdoor.com, System"
    "G107" = "backup@backdoor.com, System"
    "G108" = "backup@backdoor.com, System"
    "G113" = "admin@admin.com, dnasr281@gmail.com"
    "G132" = "backup@backdoor.com, System"
    "G133" = "backup@backdoor.com, System"
    "G134" = "backup@backdoor.com, System"
    "G139" = "admin@admin.com, dnasr281@gmail.com"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
